# Auto-generated edit script: applies the cryptos.xlsx price/volume refresh
# (and the Stacks/RenderToken row swap) described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.731.94"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").Value = "2.640.23"
$ws.Range("E3").Value = "  -1.41%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'596.08"
$ws.Range("E5").Value = "  -0.80%  "

$ws.Range("D6").Value = "'168.92"
$ws.Range("E6").Value = "  +0.76%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -2.04%  "

$ws.Range("D9").Value = "2.640.83"
$ws.Range("E9").Value = "  -1.33%  "

$ws.Range("E10").Value = "  -2.05%  "

$ws.Range("E11").Value = "  +1.30%  "

$ws.Range("E12").Value = "  +0.60%  "

$ws.Range("E13").Value = "  +0.29%  "

$ws.Range("D14").Value = "'27.76"
$ws.Range("E14").Value = "  -0.61%  "

$ws.Range("D15").Value = "3.111.59"
$ws.Range("E15").Value = "  -1.69%  "

$ws.Range("E16").Value = "  -1.31%  "

$ws.Range("D17").Value = "67.450.58"
$ws.Range("E17").Value = "  -0.16%  "

$ws.Range("D18").Value = "2.641.16"
$ws.Range("E18").Value = "  -1.55%  "

$ws.Range("E19").Value = "  +3.01%  "

$ws.Range("D20").Value = "'8.10"
$ws.Range("E20").Value = "  +2.60%  "

$ws.Range("D21").Value = "'358.77"
$ws.Range("E21").Value = "  -1.79%  "

$ws.Range("D22").Value = "'4.35"
$ws.Range("E22").Value = "  -1.02%  "

$ws.Range("E23").Value = "  -2.21%  "

$ws.Range("E24").Value = "  -3.87%  "

$ws.Range("D25").Value = "'10.41"
$ws.Range("E25").Value = "  +1.84%  "

$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("D27").Value = "'70.05"
$ws.Range("E27").Value = "  -1.22%  "

$ws.Range("D28").Value = "2.772.97"
$ws.Range("E28").Value = "  -2.02%  "

$ws.Range("E29").Value = "  +0.42%  "

$ws.Range("E30").Value = "  -1.21%  "

$ws.Range("D31").Value = "'551.33"
$ws.Range("E31").Value = "  -0.44%  "

$ws.Range("E32").Value = "  -0.29%  "

$ws.Range("E33").Value = "  -2.40%  "

$ws.Range("E34").Value = "  -1.29%  "

$ws.Range("E35").Value = "  +5.11%  "

$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("E37").Value = "  -2.73%  "

$ws.Range("D38").Value = "'157.90"
$ws.Range("E38").Value = "  +1.56%  "

$ws.Range("D39").Value = "'19.10"
$ws.Range("E39").Value = "  -2.15%  "

$ws.Range("D40").Value = "'0.368"
$ws.Range("E40").Value = "  -1.56%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "'5.26"
$ws.Range("E41").Value = "  -1.00%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.82"
$ws.Range("E42").Value = "  -0.99%  "

$ws.Range("E43").Value = "  +2.09%  "

$ws.Range("E44").Value = "  +0.09%  "

$ws.Range("E45").Value = "  -2.93%  "

$ws.Range("D46").Value = "0.0₆0304"
$ws.Range("E46").Value = "  +1.08%  "

$ws.Range("D47").Value = "'153.54"
$ws.Range("E47").Value = "  -0.24%  "

$ws.Range("E48").Value = "  -1.24%  "

$ws.Range("E49").Value = "  -1.28%  "

$ws.Range("D50").Value = "'1.72"
$ws.Range("E50").Value = "  -0.81%  "
